$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their literal text representation (the sheet stores
# prices/percentages as text, e.g. "34.94" / "1.967.55" / "  +0.33%  ", and we
# must not let Excel reinterpret numeric-looking strings as numbers, which would
# both change the stored type and silently round/alter the digits).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.574.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.37%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.94"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +6.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.302"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0698"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0955"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.071.00"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.23"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.807.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.648"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.545.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.24"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.33%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0534"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.685"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.399.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.51"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.82"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.963"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.973.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.13"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.25%  "
